# The deck ships with two theme parts:
#   ppt/theme/theme1.xml  -> bound to the (only) slide master -> "Integral" theme
#   ppt/theme/theme2.xml  -> bound to the notes master        -> "Office Theme"
#
# The target edit swaps the two themes' contents: the slide master's theme
# becomes the stock "Office Theme" palette (the notes master's theme is left
# untouched by this runtime's object model, which only exposes the deck's
# primary/slide-master theme for editing).
#
# PowerPoint's ColorScheme/ColorFormat.RGB uses the classic VBA "RGB" packed
# integer (low byte = Red, mid byte = Green, high byte = Blue), i.e. the hex
# literal must be written as 0xBBGGRR to land a given #RRGGBB value in the
# OOXML <a:srgbClr val="RRGGBB"/>.

$p = $ppt.ActivePresentation
$master = $p.Designs.Item(1).SlideMaster
$cs = $master.ColorScheme

# dk1   -> 000000
$cs.Colors(1).RGB = 0x000000
# lt1   -> FFFFFF
$cs.Colors(2).RGB = 0xFFFFFF
# dk2   -> 44546A
$cs.Colors(3).RGB = 0x6A5444
# lt2   -> E7E6E6
$cs.Colors(4).RGB = 0xE6E6E7
# accent1 -> 5B9BD5
$cs.Colors(5).RGB = 0xD59B5B
# accent2 -> ED7D31
$cs.Colors(6).RGB = 0x317DED
# accent3 -> A5A5A5
$cs.Colors(7).RGB = 0xA5A5A5
# accent4 -> FFC000
$cs.Colors(8).RGB = 0x00C0FF
# accent5 -> 4472C4
$cs.Colors(9).RGB = 0xC47244
# accent6 -> 70AD47
$cs.Colors(10).RGB = 0x47AD70
# hlink    -> 0563C1
$cs.Colors(11).RGB = 0xC16305
# folHlink -> 954F72
$cs.Colors(12).RGB = 0x724F95
